# Auto-generated edit script: updates crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.976.77"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "3.173.16"
$ws.Range("E3").Value = "  -4.43%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.36"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.41"
$ws.Range("E6").Value = "  -3.98%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.170.01"
$ws.Range("E8").Value = "  -4.47%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -4.49%  "
$ws.Range("E11").Value = "  -3.46%  "
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.85"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "3.694.83"
$ws.Range("E15").Value = "  -4.45%  "
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "3.174.65"
$ws.Range("E17").Value = "  -4.40%  "
$ws.Range("D18").Value = "62.910.17"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.63"
$ws.Range("E19").Value = "  -3.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.21"
$ws.Range("E20").Value = "  -3.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.94"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("E22").Value = "  -3.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.65"
$ws.Range("E23").Value = "  -6.16%  "
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -2.95%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  -4.48%  "
$ws.Range("E30").Value = "  -6.26%  "
$ws.Range("E31").Value = "  -6.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.30"
$ws.Range("E32").Value = "  -5.76%  "
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("E35").Value = "  -6.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.84"
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.12"
$ws.Range("E37").Value = "  -3.16%  "
$ws.Range("D38").Value = "0.0₃0712"
$ws.Range("E38").Value = "  -4.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0389"
$ws.Range("E39").Value = "  -2.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "405.55"
$ws.Range("E40").Value = "  -6.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.69"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("D44").Value = "2.789.28"
$ws.Range("E44").Value = "  -9.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.253"
$ws.Range("E45").Value = "  -3.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.14"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.79"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.70"
$ws.Range("E49").Value = "  -5.62%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.81"
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("E51").Value = "  -2.22%  "
